{"js": "// Lesson 4 Understanding JSX | Complete React Course in Hindi\n// Insert six new notes paragraphs right after the paragraph that ends with\n// \"used to avoid bugs in our scripts\" and before the trailing empty\n// paragraph at the end of the document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its (unique) text.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"used to avoid bugs in our scripts\") !== -1) {\n    anchor = p;\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph ending with 'used to avoid bugs in our scripts'\");\n}\n\nconst newParagraphTexts = [\n  \"Use htmlFor instead of for.\",\n  \"Use className instead of class.\",\n  \"Use tabIndex instead of tabindex\",\n  \"It is mandatory that tag or componeent should be closed in jsx even if it is an empty tag by /\",\n  \"React use webpack for liveserver \",\n  \"And use babel for compilation of component to js\"\n];\n\n// Insert paragraphs in order, each after the previous one, so the final\n// order matches the list above.\nlet current = anchor;\nfor (const text of newParagraphTexts) {\n  current = current.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Lesson 4 Understanding JSX | Complete React Course in Hindi\n# Insert six new notes paragraphs right after the paragraph that ends with\n# \"used to avoid bugs in our scripts\" and before the trailing empty\n# paragraph at the end of the document body.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its text instead of a hard-coded index so\n# the script is resilient to any earlier structural differences.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*used to avoid bugs in our scripts*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph ending with 'used to avoid bugs in our scripts'\"\n}\n\n$newParagraphTexts = @(\n    \"Use htmlFor instead of for.\",\n    \"Use className instead of class.\",\n    \"Use tabIndex instead of tabindex\",\n    \"It is mandatory that tag or componeent should be closed in jsx even if it is an empty tag by /\",\n    \"React use webpack for liveserver \",\n    \"And use babel for compilation of component to js\"\n)\n\n$currentIndex = $anchorIndex\nforeach ($text in $newParagraphTexts) {\n    $d.Paragraphs.Item($currentIndex).Range.InsertParagraphAfter()\n    $currentIndex = $currentIndex + 1\n    $d.Paragraphs.Item($currentIndex).Range.Text = $text\n}\n"}
